# The "Export" sheet lists account balances (Conta/Nome/Saldo) sorted by
# descending Saldo. The row for account 000772433 (MARCELO), currently at
# worksheet row 15 with Saldo 7707.44, needs to move down past the next 10
# rows (CESAR ... LAURA) to sit right before EVANGELINA (currently row 26),
# and its Saldo needs to become 707.44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the row just above EVANGELINA (row 26) by inserting a
#    blank row there, shifting EVANGELINA and everything below it down by one.
#    Using an explicit A:C range (rather than the whole row) keeps the insert
#    limited to the columns that are actually in use.
$ws.Range("A26:C26").Insert()

# 2) Cut the MARCELO row (still at row 15, since nothing above row 26 moved
#    yet) and paste it into the newly created blank row.
$src = $ws.Range("A15:C15")
$dst = $ws.Range("A26:C26")
$src.Cut($dst)

# 3) Remove the now-empty row 15 left behind by the cut, shifting rows
#    16-26 back up by one. MARCELO's row now lands on row 25.
$ws.Range("A15:C15").Delete()

# 4) Update MARCELO's Saldo value per the edit (7707.44 -> 707.44).
$ws.Range("C25").Value2 = 707.44
